$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.060112
$ws.Range("H2").Value = 0.180336
$ws.Range("I2").Value = 0.02822008145728641
$ws.Range("J2").Value = 0.02822008145728641
$ws.Range("M2").Value = 1.599392
$ws.Range("N2").Value = 4.798176
$ws.Range("O2").Value = 0.03952976301548796
$ws.Range("P2").Value = 0.03952976301548796
$ws.Range("Q2").Value = 0.096142651904
$ws.Range("R2").Value = 0.8652838671359999
$ws.Range("S2").Value = 0.001115533132284298
$ws.Range("T2").Value = 0.001115533132284298

$ws.Range("G3").Value = 0.060112
$ws.Range("H3").Value = 0.180336
$ws.Range("I3").Value = 0.02822008145728641
$ws.Range("J3").Value = 0.02822008145728641
$ws.Range("O3").Value = 0.4638329693976876
$ws.Range("P3").Value = 0.4638329693976876
$ws.Range("Q3").Value = 1.128115331754667
$ws.Range("R3").Value = 10.153037985792
$ws.Range("S3").Value = 0.01308940417897778
$ws.Range("T3").Value = 0.01308940417897778

$ws.Range("G4").Value = 0.060112
$ws.Range("H4").Value = 0.180336
$ws.Range("I4").Value = 0.02822008145728641
$ws.Range("J4").Value = 0.02822008145728641
$ws.Range("M4").Value = 20.09416733333333
$ws.Range("N4").Value = 60.28250199999999
$ws.Range("O4").Value = 0.4966372675868244
$ws.Range("P4").Value = 0.4966372675868245
$ws.Range("Q4").Value = 1.207900586741333
$ws.Range("R4").Value = 10.871105280672
$ws.Range("S4").Value = 0.01401514414602433
$ws.Range("T4").Value = 0.01401514414602434

$ws.Range("I5").Value = 0.716133390648984
$ws.Range("J5").Value = 0.716133390648984
$ws.Range("M5").Value = 1.599392
$ws.Range("N5").Value = 4.798176
$ws.Range("O5").Value = 0.03952976301548796
$ws.Range("P5").Value = 0.03952976301548796
$ws.Range("Q5").Value = 2.439786128832
$ws.Range("R5").Value = 21.958075159488
$ws.Range("S5").Value = 0.0283085832198322
$ws.Range("T5").Value = 0.0283085832198322

$ws.Range("I6").Value = 0.716133390648984
$ws.Range("J6").Value = 0.716133390648984
$ws.Range("O6").Value = 0.4638329693976876
$ws.Range("P6").Value = 0.4638329693976876
$ws.Range("S6").Value = 0.3321662770695524
$ws.Range("T6").Value = 0.3321662770695524

$ws.Range("I7").Value = 0.716133390648984
$ws.Range("J7").Value = 0.716133390648984
$ws.Range("M7").Value = 20.09416733333333
$ws.Range("N7").Value = 60.28250199999999
$ws.Range("O7").Value = 0.4966372675868244
$ws.Range("P7").Value = 0.4966372675868245
$ws.Range("Q7").Value = 30.652567181964
$ws.Range("R7").Value = 275.873104637676
$ws.Range("S7").Value = 0.3556585303595993
$ws.Range("T7").Value = 0.3556585303595993

$ws.Range("G8").Value = 0.5445563333333333
$ws.Range("H8").Value = 1.633669
$ws.Range("I8").Value = 0.2556465278937297
$ws.Range("J8").Value = 0.2556465278937297
$ws.Range("M8").Value = 1.599392
$ws.Range("N8").Value = 4.798176
$ws.Range("O8").Value = 0.03952976301548796
$ws.Range("P8").Value = 0.03952976301548796
$ws.Range("Q8").Value = 0.8709590430826666
$ws.Range("R8").Value = 7.838631387744
$ws.Range("S8").Value = 0.01010564666337147
$ws.Range("T8").Value = 0.01010564666337147

$ws.Range("G9").Value = 0.5445563333333333
$ws.Range("H9").Value = 1.633669
$ws.Range("I9").Value = 0.2556465278937297
$ws.Range("J9").Value = 0.2556465278937297
$ws.Range("O9").Value = 0.4638329693976876
$ws.Range("P9").Value = 0.4638329693976876
$ws.Range("Q9").Value = 10.21962916950756
$ws.Range("R9").Value = 91.97666252556802
$ws.Range("S9").Value = 0.1185772881491574
$ws.Range("T9").Value = 0.1185772881491574

$ws.Range("G10").Value = 0.5445563333333333
$ws.Range("H10").Value = 1.633669
$ws.Range("I10").Value = 0.2556465278937297
$ws.Range("J10").Value = 0.2556465278937297
$ws.Range("M10").Value = 20.09416733333333
$ws.Range("N10").Value = 60.28250199999999
$ws.Range("O10").Value = 0.4966372675868244
$ws.Range("P10").Value = 0.4966372675868245
$ws.Range("Q10").Value = 10.94240608442644
$ws.Range("R10").Value = 98.48165475983799
$ws.Range("S10").Value = 0.1269635930812008
$ws.Range("T10").Value = 0.1269635930812008
